$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 109, shifting the existing
# rows 109-114 down to 110-115 (the data table is sorted by date, and a
# newer weekly reading is being added in chronological order).
$ws.Rows.Item(109).Insert()

# Populate the newly inserted row 109 with the new weekly reading.
$ws.Cells.Item(109, 1).Value = 5
$ws.Cells.Item(109, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(109, 3).Value = "Maule"
$ws.Cells.Item(109, 4).Value = 44610
$ws.Cells.Item(109, 5).Value = 7
$ws.Cells.Item(109, 6).Value = 100112030
$ws.Cells.Item(109, 7).Value = "Poroto granado"
$ws.Cells.Item(109, 8).Value = "Sin especificar"
$ws.Cells.Item(109, 9).Value = "Primera"
$ws.Cells.Item(109, 10).Value = 400
$ws.Cells.Item(109, 11).Value = 17000
$ws.Cells.Item(109, 12).Value = 17000
$ws.Cells.Item(109, 13).Value = 17000
$ws.Cells.Item(109, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(109, 15).Value = "Región del Maule"
$ws.Cells.Item(109, 16).Value = 680
$ws.Cells.Item(109, 17).Value = 25
$ws.Cells.Item(109, 18).Value = "Hortaliza"
